$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.110.85"
$ws.Range("E2").Value = "  -1.34%  "

$ws.Range("D3").Value = "1.992.72"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4972"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4192"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.67"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08908"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.27"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").Value = "1.997.15"
$ws.Range("E13").Value = "  -1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.011"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.437"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.015"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.53%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.41"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.53%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001107"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06765"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.59"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.52%  "

$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.973"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").Value = "29.115.68"
$ws.Range("E23").Value = "  -1.43%  "

$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.290"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("D26").Value = "2.229.78"
$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.85"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.18"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.290"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.260"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.39"
$ws.Range("D31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.047"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09850"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.533"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.829"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.743"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02423"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.323"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.070"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06398"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6481"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.49"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1988"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.011"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6213"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.351"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.46"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.186"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.99%  "

$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000336"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.65%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.152"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.64%  "

